$wb = $excel.ActiveWorkbook

# Switch to the "Service" worksheet and delete row 8 (the extra_service_id /
# "Foregin Key (ExtraService->extra_id)" field row) — the whole row shifts
# up, taking the rows below it (total_payment, comments, pat_at_home, ...)
# with it.
$ws = $wb.Worksheets.Item("Service")
$ws.Activate()
$ws.Rows("8").Delete()

# Leave the selection on C8 (where "comments"/total_payment's attribute cell
# now sits) and make sure this sheet is the one showing as active when the
# file is reopened.
$ws.Range("C8").Select()
